$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Excel_vs_ML")
$ws1.Range("H2").Value = 100271.72
$ws1.Range("L2").Value = 99.92
$ws1.Range("M2").Value = 78.65000000000001
$ws1.Range("O2").Value = "On Track"
$ws1.Range("H3").Value = 445313.96
$ws1.Range("L3").Value = 100.28
$ws1.Range("M3").Value = 107280.67
$ws1.Range("N3").Value = 9752.790000000001
$ws1.Range("O3").Value = "On Track"
$ws1.Range("P3").Value = "On Track"
$ws1.Range("Q3").Value = 0
$ws1.Range("H4").Value = 294599.18
$ws1.Range("L4").Value = 100.02
$ws1.Range("M4").Value = -58.8
$ws1.Range("O4").Value = "On Track"
$ws1.Range("H5").Value = 503642.69
$ws1.Range("L5").Value = 98.93000000000001
$ws1.Range("M5").Value = 5454.68
$ws1.Range("O5").Value = "On Track"
$ws1.Range("H6").Value = 373866.06
$ws1.Range("L6").Value = 97.77
$ws1.Range("M6").Value = 8543.190000000001
$ws1.Range("H7").Value = 511982.16
$ws1.Range("L7").Value = 97.39
$ws1.Range("M7").Value = 13713.78
$ws1.Range("O7").Value = "On Track"
$ws1.Range("H8").Value = 63085.98
$ws1.Range("L8").Value = 160
$ws1.Range("M8").Value = 73754.94
$ws1.Range("N8").Value = 1756.07
$ws1.Range("P8").Value = "On Track"
$ws1.Range("Q8").Value = 0
$ws1.Range("H10").Value = 110328
$ws1.Range("L10").Value = 100.07
$ws1.Range("M10").Value = 338530.15
$ws1.Range("N10").Value = 7872.79
$ws1.Range("O10").Value = "On Track"
$ws1.Range("P10").Value = "On Track"
$ws1.Range("H12").Value = 250343.25
$ws1.Range("L12").Value = 99.11
$ws1.Range("M12").Value = 309463.48
$ws1.Range("N12").Value = 6876.97
$ws1.Range("O12").Value = "On Track"
$ws1.Range("P12").Value = "On Track"
$ws1.Range("H13").Value = 51709.52
$ws1.Range("L13").Value = 107.52
$ws1.Range("M13").Value = 117668.02
$ws1.Range("N13").Value = 2028.76
$ws1.Range("P13").Value = "On Track"
$ws1.Range("Q13").Value = 0
$ws1.Range("H14").Value = 335066.17
$ws1.Range("L14").Value = 100.04
$ws1.Range("M14").Value = -130.08
$ws1.Range("O14").Value = "On Track"
$ws1.Range("H16").Value = 158329.38
$ws1.Range("L16").Value = 97.97
$ws1.Range("M16").Value = 3275.72
$ws1.Range("H17").Value = 96927.23
$ws1.Range("L17").Value = 101.41
$ws1.Range("M17").Value = -1343.42
$ws1.Range("O17").Value = "On Track"
$ws1.Range("H20").Value = 239614.62
$ws1.Range("L20").Value = 133.95
$ws1.Range("M20").Value = 184065.29
$ws1.Range("N20").Value = 3539.72
$ws1.Range("O20").Value = "Overpacing"
$ws1.Range("P20").Value = "On Track"
$ws1.Range("H23").Value = 106044.44
$ws1.Range("L23").Value = 100.09
$ws1.Range("M23").Value = -93.37
$ws1.Range("O23").Value = "On Track"
$ws1.Range("H24").Value = 274327.61
$ws1.Range("L24").Value = 97.89
$ws1.Range("M24").Value = 5920.98
$ws1.Range("O24").Value = "On Track"
$ws1.Range("H25").Value = 287064.14
$ws1.Range("L25").Value = 100.42
$ws1.Range("M25").Value = -1200.91
$ws1.Range("O25").Value = "On Track"
$ws1.Range("H26").Value = 235916.39
$ws1.Range("L26").Value = 106.67
$ws1.Range("M26").Value = -14745.33
$ws1.Range("H27").Value = 233908.29
$ws1.Range("L27").Value = 65.88
$ws1.Range("M27").Value = 343900.59
$ws1.Range("N27").Value = 10746.89
$ws1.Range("O27").Value = "Underpacing"
$ws1.Range("P27").Value = "On Track"
$ws1.Range("Q27").Value = 0
$ws1.Range("H28").Value = 198664.26
$ws1.Range("L28").Value = 100.96
$ws1.Range("M28").Value = -1886.9
$ws1.Range("O28").Value = "On Track"
$ws1.Range("H29").Value = 162193.44
$ws1.Range("L29").Value = 95
$ws1.Range("M29").Value = 8536.57
$ws1.Range("O29").Value = "Underpacing"
$ws1.Range("H30").Value = 348450.7
$ws1.Range("L30").Value = 97.97
$ws1.Range("M30").Value = 7209.56
$ws1.Range("O30").Value = "On Track"
$ws1.Range("H31").Value = 406291.6
$ws1.Range("L31").Value = 97
$ws1.Range("M31").Value = 12565.79
$ws1.Range("O31").Value = "On Track"
$ws1.Range("H32").Value = 366396.75
$ws1.Range("L32").Value = 121.43
$ws1.Range("M32").Value = 208343.52
$ws1.Range("N32").Value = 3655.15
$ws1.Range("O32").Value = "Overpacing"
$ws1.Range("P32").Value = "On Track"
$ws1.Range("H33").Value = 66136.17
$ws1.Range("L33").Value = 102.25
$ws1.Range("M33").Value = 63220.8
$ws1.Range("N33").Value = 1915.78
$ws1.Range("O33").Value = "On Track"
$ws1.Range("P33").Value = "On Track"
$ws1.Range("H34").Value = 553167.48
$ws1.Range("L34").Value = 93.45
$ws1.Range("M34").Value = 38783.01
$ws1.Range("H35").Value = 383062.2
$ws1.Range("L35").Value = 98.13
$ws1.Range("M35").Value = 7319.61
$ws1.Range("O35").Value = "On Track"
$ws1.Range("H36").Value = 124810.24
$ws1.Range("L36").Value = 97.3
$ws1.Range("M36").Value = 3458.55
$ws1.Range("O36").Value = "On Track"
$ws1.Range("H38").Value = 422599.07
$ws1.Range("L38").Value = 108.03
$ws1.Range("M38").Value = -31402.5
$ws1.Range("O38").Value = "Overpacing"
$ws1.Range("H39").Value = 593681.45
$ws1.Range("L39").Value = 99.55
$ws1.Range("M39").Value = 2697.48
$ws1.Range("O39").Value = "On Track"
$ws1.Range("H40").Value = 43416.78
$ws1.Range("L40").Value = 61.84
$ws1.Range("M40").Value = 76672.67
$ws1.Range("N40").Value = 2839.73
$ws1.Range("P40").Value = "On Track"
$ws1.Range("H43").Value = 334645.05
$ws1.Range("L43").Value = 95
$ws1.Range("M43").Value = 17612.41
$ws1.Range("O43").Value = "On Track"
$ws1.Range("H45").Value = 196607.05
$ws1.Range("L45").Value = 76.84
$ws1.Range("M45").Value = 158003.96
$ws1.Range("N45").Value = 7182
$ws1.Range("O45").Value = "Underpacing"
$ws1.Range("P45").Value = "On Track"
$ws1.Range("Q45").Value = 0

$ws2 = $wb.Worksheets.Item("Feature_Importance")
$ws2.Range("B2").Value = 0.2159633534796498
$ws2.Range("A3").Value = "Spend_Velocity"
$ws2.Range("B3").Value = 0.1603340188252249
$ws2.Range("A4").Value = "Flight_Days"
$ws2.Range("B4").Value = 0.1550880643428183
$ws2.Range("B5").Value = 0.1407617252684994
$ws2.Range("A6").Value = "Days_Elapsed"
$ws2.Range("B6").Value = 0.1337881870750038
$ws2.Range("A7").Value = "Spend_to_Date"
$ws2.Range("B7").Value = 0.1271158795313068
$ws2.Range("B8").Value = 0.06694877147749694

$ws3 = $wb.Worksheets.Item("Exec_Summary")
$ws3.Range("B2").Value = 0.833
$ws3.Range("B3").Value = 11
$ws3.Range("B4").Value = 0
$ws3.Range("B6").Value = 0
